$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values must stay as text
# (mirrors the source data which is stored as text, not numbers)
$textCells = @("D4", "D5", "D6", "D10", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D26", "D27", "D28", "D31", "D34", "D38", "D39", "D40", "D41", "D43", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.870.09'
$ws.Range("E2").Value = '  +2.80%  '
$ws.Range("D3").Value = '2.531.99'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '593.37'
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").Value = '176.55'
$ws.Range("E6").Value = '  +5.84%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").Value = '2.531.43'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +2.28%  '
$ws.Range("E11").Value = '  +2.38%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '0.345'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = '26.89'
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '2.991.35'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = '67.686.63'
$ws.Range("E17").Value = '  +2.82%  '
$ws.Range("D18").Value = '2.530.17'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '8.04'
$ws.Range("E19").Value = '  +5.55%  '
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("D21").Value = '361.14'
$ws.Range("E21").Value = '  +4.39%  '
$ws.Range("D22").Value = '4.21'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").Value = '4.66'
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("E24").Value = '  +3.31%  '
$ws.Range("D26").Value = '10.33'
$ws.Range("E26").Value = '  +3.65%  '
$ws.Range("D27").Value = '70.97'
$ws.Range("E27").Value = '  +3.30%  '
$ws.Range("D28").Value = '0.996'
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D30").Value = '0.0₃0991'
$ws.Range("E30").Value = '  +2.04%  '
$ws.Range("D31").Value = '554.49'
$ws.Range("E31").Value = '  +5.87%  '
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("E33").Value = '  +3.46%  '
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +2.96%  '
$ws.Range("E35").Value = '  +0.88%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  +2.47%  '
$ws.Range("D38").Value = '155.81'
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("D39").Value = '18.79'
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("D40").Value = '18.60'
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("D41").Value = '0.357'
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("E42").Value = '  +3.62%  '
$ws.Range("D43").Value = '5.19'
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("E44").Value = '  +5.60%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '0.563'
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '147.45'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = '0.0₆0281'
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("D49").Value = '3.73'
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '0.0759'
$ws.Range("E51").Value = '  +0.47%  '
